# The "LINK & IDEA" reference slide (a leftover dev note pointing at an
# animation-library link, unrelated to the CrazyMovies wireframe deck) is
# removed from the presentation. It was the last slide (slide 10 /
# sldId 265).
$p = $ppt.ActivePresentation

$lastIndex = $p.Slides.Count
$slide = $p.Slides.Item($lastIndex)

# Delete it directly if it's the expected trailing "LINK & IDEA" slide;
# otherwise fall back to scanning for it by its title text so the script
# stays correct even if slide ordering ever shifts.
$title = ""
try {
    $title = $slide.Shapes.Item(1).TextFrame.TextRange.Text
} catch {
    $title = ""
}

if ($title -like "*LINK*IDEA*") {
    $slide.Delete()
} else {
    for ($i = $p.Slides.Count; $i -ge 1; $i--) {
        $candidate = $p.Slides.Item($i)
        $candidateTitle = ""
        try {
            $candidateTitle = $candidate.Shapes.Item(1).TextFrame.TextRange.Text
        } catch {
            $candidateTitle = ""
        }
        if ($candidateTitle -like "*LINK*IDEA*") {
            $candidate.Delete()
            break
        }
    }
}
